$d = $word.ActiveDocument

# Find the "ADS_Portfolio_Reflection.docx" list-item paragraph that we will
# insert the new entries after.
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "ADS_Portfolio_Reflection.docx") {
        $anchor = $p
    }
}

$newItems = @(
    "ADS_Portfolio_Reflection_Deck.pptx",
    "Portfolio_reflection_v4.mp4",
    "LICENSE",
    "README"
)

$r = $anchor.Range
foreach ($text in $newItems) {
    $r.InsertParagraphAfter()
    $r = $r.Next(4, 1)
    $r.Text = $text
}
